# Update attendance ("想去人数") and ticket price ("最低票价") figures
# across the four worksheets, matching the upstream data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("G3").Value  = 109
$ws1.Range("F6").Value  = 984
$ws1.Range("F8").Value  = 67
$ws1.Range("F10").Value = 920
$ws1.Range("F11").Value = 337
$ws1.Range("F14").Value = 1392
$ws1.Range("F17").Value = 2959
$ws1.Range("F18").Value = 359
$ws1.Range("F19").Value = 1582
$ws1.Range("F23").Value = 1319
$ws1.Range("F24").Value = 247
$ws1.Range("F27").Value = 383
$ws1.Range("F28").Value = 3375

# --- Sheet 2: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F6").Value  = 16
$ws2.Range("F7").Value  = 4
$ws2.Range("F8").Value  = 42
$ws2.Range("F9").Value  = 13
$ws2.Range("F13").Value = 60

# --- Sheet 3: 本地生活 (Local life) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 782

# --- Sheet 4: 全部类型 (All types, aggregated view) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value  = 782
$ws4.Range("G4").Value  = 109
$ws4.Range("F8").Value  = 984
$ws4.Range("F11").Value = 67
$ws4.Range("F15").Value = 16
$ws4.Range("F16").Value = 4
$ws4.Range("F17").Value = 42
$ws4.Range("F19").Value = 13
$ws4.Range("F21").Value = 920
$ws4.Range("F22").Value = 337
$ws4.Range("F25").Value = 1392
$ws4.Range("F28").Value = 2959
$ws4.Range("F29").Value = 359
$ws4.Range("F30").Value = 1582
$ws4.Range("F34").Value = 1319
$ws4.Range("F35").Value = 247
$ws4.Range("F40").Value = 383
$ws4.Range("F41").Value = 3375
$ws4.Range("F45").Value = 60
